$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1297.1428
$ws.Range("I19").Value = 1002
$ws.Range("K19").Value = 1002
$ws.Range("M19").Value = -827
$ws.Range("H32").Value = 4198.4
$ws.Range("H44").Value = 39999
$ws.Range("J44").Value = 39999
$ws.Range("L44").Value = 39999
$ws.Range("N44").Value = -40923
$ws.Range("H51").Value = 4999.5
$ws.Range("J51").Value = 4999
$ws.Range("L51").Value = 4999
$ws.Range("N51").Value = -5967
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H100").Value = 2643.8
$ws.Range("J100").Value = 3150
$ws.Range("L100").Value = 3150
$ws.Range("N100").Value = -4232
$ws.Range("H116").Value = 5378.0356
$ws.Range("I116").Value = 5195.56
$ws.Range("J116").Value = 6898.6665
$ws.Range("K116").Value = 5195.56
$ws.Range("L116").Value = 6898.6665
$ws.Range("M116").Value = -1753.56
$ws.Range("N116").Value = -13782.6665
$ws.Range("H135").Value = 50000388
$ws.Range("I135").Value = 55555964
$ws.Range("K135").Value = 500003676
$ws.Range("M135").Value = -500001141
$ws.Range("H137").Value = 2452.682
$ws.Range("I137").Value = 2098.6667
$ws.Range("K137").Value = 6296.000100000001
$ws.Range("M137").Value = -3746.000100000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19231
$ws.Range("H74").Value = 27029876
$ws.Range("I74").Value = 27780692
$ws.Range("K74").Value = 27780692
$ws.Range("M74").Value = -27779818
$ws.Range("H75").Value = 110000
$ws.Range("J75").Value = 110000
$ws.Range("L75").Value = 110000
$ws.Range("N75").Value = -111748
$ws.Range("H77").Value = 27029876
$ws.Range("I77").Value = 27780692
$ws.Range("K77").Value = 138903460
$ws.Range("M77").Value = -138899092
$ws.Range("H78").Value = 110000
$ws.Range("J78").Value = 110000
$ws.Range("L78").Value = 330000
$ws.Range("N78").Value = -338736
$ws.Range("H122").Value = 3884.3333
$ws.Range("I122").Value = 2951.1738
$ws.Range("K122").Value = 8853.5214
$ws.Range("M122").Value = -6403.5214
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3449.4
$ws.Range("I20").Value = 3415
$ws.Range("J20").Value = 3518.2
$ws.Range("K20").Value = 3415
$ws.Range("L20").Value = 3518.2
$ws.Range("M20").Value = -3168
$ws.Range("N20").Value = -4012.2
$ws.Range("H86").Value = 2526.5557
$ws.Range("I86").Value = 2768.647
$ws.Range("J86").Value = 2115
$ws.Range("K86").Value = 2768.647
$ws.Range("L86").Value = 2115
$ws.Range("M86").Value = -1645.647
$ws.Range("N86").Value = -4361
$ws.Range("H89").Value = 2526.5557
$ws.Range("I89").Value = 2768.647
$ws.Range("J89").Value = 2115
$ws.Range("K89").Value = 13843.235
$ws.Range("L89").Value = 10575
$ws.Range("M89").Value = -8227.235000000001
$ws.Range("N89").Value = -21807
$ws.Range("H99").Value = 1510.3928
$ws.Range("I99").Value = 1468.8334
$ws.Range("K99").Value = 1468.8334
$ws.Range("M99").Value = 29.16660000000002
$ws.Range("H134").Value = 31251148
$ws.Range("J134").Value = 997
$ws.Range("L134").Value = 2991
$ws.Range("N134").Value = -8061
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 845
$ws.Range("I22").Value = 819.1667
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 819.1667
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -469.1667
$ws.Range("N22").Value = -1700
$ws.Range("H99").Value = 16136
$ws.Range("I99").Value = 17442.715
$ws.Range("K99").Value = 17442.715
$ws.Range("M99").Value = -15944.715
$ws.Range("H126").Value = 16136
$ws.Range("I126").Value = 17442.715
$ws.Range("K126").Value = 52328.145
$ws.Range("M126").Value = -49858.145
$ws.Range("H134").Value = 10917682
$ws.Range("I134").Value = 11413622
$ws.Range("K134").Value = 34240866
$ws.Range("M134").Value = -34238331
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 114749.25
$ws.Range("J37").Value = 114749.25
$ws.Range("L37").Value = 344247.75
$ws.Range("N37").Value = -344471.75
$ws.Range("H39").Value = 726
$ws.Range("I39").Value = 726
$ws.Range("K39").Value = 2178
$ws.Range("M39").Value = -1884
$ws.Range("H140").Value = 2430
$ws.Range("I140").Value = 2430
$ws.Range("K140").Value = 7290
$ws.Range("M140").Value = -2110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2499.5
$ws.Range("I80").Value = 2499.5
$ws.Range("K80").Value = 2499.5
$ws.Range("M80").Value = -1501.5
$ws.Range("H83").Value = 2499.5
$ws.Range("I83").Value = 2499.5
$ws.Range("K83").Value = 12497.5
$ws.Range("M83").Value = -7505.5
$ws.Range("H132").Value = 4035335.8
$ws.Range("I132").Value = 5002848.5
$ws.Range("J132").Value = 4031.3333
$ws.Range("K132").Value = 15008545.5
$ws.Range("L132").Value = 12093.9999
$ws.Range("M132").Value = -15006015.5
$ws.Range("N132").Value = -17153.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3353.889
$ws.Range("I22").Value = 3399.375
$ws.Range("J22").Value = 2990
$ws.Range("K22").Value = 3399.375
$ws.Range("L22").Value = 2990
$ws.Range("M22").Value = -3104.375
$ws.Range("N22").Value = -3580
$ws.Range("H27").Value = 3353.889
$ws.Range("I27").Value = 3399.375
$ws.Range("J27").Value = 2990
$ws.Range("K27").Value = 3399.375
$ws.Range("L27").Value = 2990
$ws.Range("M27").Value = -3292.375
$ws.Range("N27").Value = -3204
$ws.Range("H46").Value = 2344.4443
$ws.Range("I46").Value = 2375
$ws.Range("K46").Value = 2375
$ws.Range("M46").Value = -2187
$ws.Range("H55").Value = 152.8
$ws.Range("J55").Value = 104
$ws.Range("L55").Value = 104
$ws.Range("N55").Value = -450
$ws.Range("H61").Value = 6348
$ws.Range("I61").Value = 6348
$ws.Range("K61").Value = 6348
$ws.Range("M61").Value = -6146
$ws.Range("H68").Value = 8775363
$ws.Range("I68").Value = 26315788
$ws.Range("K68").Value = 26315788
$ws.Range("M68").Value = -26315039
$ws.Range("H71").Value = 8775363
$ws.Range("I71").Value = 26315788
$ws.Range("K71").Value = 131578940
$ws.Range("M71").Value = -131575196
$ws.Range("H82").Value = 876.88464
$ws.Range("I82").Value = 1003.4
$ws.Range("J82").Value = 704.36365
$ws.Range("K82").Value = 1003.4
$ws.Range("L82").Value = 704.36365
$ws.Range("M82").Value = -642.4
$ws.Range("N82").Value = -1426.36365
$ws.Range("H85").Value = 876.88464
$ws.Range("I85").Value = 1003.4
$ws.Range("J85").Value = 704.36365
$ws.Range("K85").Value = 1003.4
$ws.Range("L85").Value = 704.36365
$ws.Range("M85").Value = 244.6
$ws.Range("N85").Value = -3200.36365
$ws.Range("H93").Value = 1844.8889
$ws.Range("I93").Value = 1109
$ws.Range("K93").Value = 1109
$ws.Range("M93").Value = 139
$ws.Range("H113").Value = 6348
$ws.Range("I113").Value = 6348
$ws.Range("K113").Value = 6348
$ws.Range("M113").Value = -4178
$ws.Range("H132").Value = 16672654
$ws.Range("I132").Value = 17863272
$ws.Range("K132").Value = 53589816
$ws.Range("M132").Value = -53587286
$ws.Range("H139").Value = 298994.2
$ws.Range("J139").Value = 298992.75
$ws.Range("L139").Value = 298992.75
$ws.Range("N139").Value = -309272.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2700
$ws.Range("J8").Value = 2700
$ws.Range("L8").Value = 2700
$ws.Range("N8").Value = -2980
$ws.Range("H113").Value = 945.73334
$ws.Range("I113").Value = 799.0714
$ws.Range("K113").Value = 2397.2142
$ws.Range("M113").Value = -227.2142000000003
$ws.Range("H122").Value = 2188.4783
$ws.Range("I122").Value = 1784.7059
$ws.Range("J122").Value = 3332.5
$ws.Range("K122").Value = 5354.1177
$ws.Range("L122").Value = 9997.5
$ws.Range("M122").Value = -2904.1177
$ws.Range("N122").Value = -14897.5
$ws.Range("H132").Value = 15155295
$ws.Range("I132").Value = 22729262
$ws.Range("K132").Value = 68187786
$ws.Range("M132").Value = -68185256

Write-Output "Applied all cell updates"